$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 10 with values, matching the style of the existing "A" column (bold/centered/bordered)
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 49
$ws.Range("C10").Value = 51
$ws.Range("D10").Value = 53
$ws.Range("E10").Value = 54
$ws.Range("F10").Value = 52
$ws.Range("G10").Value = 50

# Copy the style used by the other row-index cells (A2:A9) onto A10
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = 9
